# Regenerate merged AHB files
# 1. Rename the diff-comparison header columns:
#      *_old  -> *_FV2210
#      *_new  -> *_FV2304
#    (column K / "diff" stays as-is)
# 2. Wrap the used range in an Excel Table ("Table1").
# 3. Freeze the header row (row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2210", "Segmentgruppe_FV2210", "Segment_FV2210", "Datenelement_FV2210", "Segment ID_FV2210",
    "Code_FV2210", "Qualifier_FV2210", "Beschreibung_FV2210", "Bedingungsausdruck_FV2210", "Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304", "Segmentgruppe_FV2304", "Segment_FV2304", "Datenelement_FV2304", "Segment ID_FV2304",
    "Code_FV2304", "Qualifier_FV2304", "Beschreibung_FV2304", "Bedingungsausdruck_FV2304", "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn the data range into a proper Excel Table (ListObject) named "Table1"
$tableRange = $ws.Range("A1:U70")
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"
$table.TableStyle = ""

# Freeze the first (header) row
$ws.Range("A2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
